$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 8, pushing existing rows 8-115 down to 9-116
$ws.Rows("8:8").Insert()

# Populate the newly inserted row 8 with the new weekly record
$ws.Range("A8").Value = 7
$ws.Range("B8").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C8").Value = "Ñuble"
$ws.Range("D8").Value = 44963
$ws.Range("D8").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E8").Value = 16
$ws.Range("F8").Value = 100112030
$ws.Range("G8").Value = "Poroto granado"
$ws.Range("H8").Value = "Sin especificar"
$ws.Range("I8").Value = "Primera"
$ws.Range("J8").Value = 50
$ws.Range("K8").Value = 35000
$ws.Range("L8").Value = 35000
$ws.Range("M8").Value = 35000
$ws.Range("N8").Value = "$/saco 25 kilos"
$ws.Range("O8").Value = "Región del Maule"
$ws.Range("P8").Value = 1400
$ws.Range("Q8").Value = 25
$ws.Range("R8").Value = "Hortaliza"
